$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Nur"
$ws.Range("B5").Value = 24

$ws.Range("A6").Value = "Biplob"
$ws.Range("B6").Value = 31

$ws.Range("B6").Select()
